$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet currently ends with a "Total" row at row 16 (SUM(D2:D15) / SUM(F2:F15)).
# Two new timesheet entries need to be inserted above it (new rows 16 & 17),
# pushing the Total row down to row 18 and widening its SUM ranges to D2:D17 / F2:F17.

# Insert two blank rows at row 16 (this pushes the old Total row from 16 -> 18).
$ws.Rows.Item(16).Insert()
$ws.Rows.Item(16).Insert()

# Give the two new rows the same cell formatting as the last data row (row 15)
# so they pick up the correct date / time / hours / currency styles.
$ws.Range("A15:F15").Copy()
$ws.Range("A16:F17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New data row 16: 2023-12-26, 15:00 -> 20:00, rate 10
$ws.Cells.Item(16,1).Value = 45286
$ws.Cells.Item(16,2).Value = 0.625
$ws.Cells.Item(16,3).Value = 0.833333333333333
$ws.Cells.Item(16,4).Formula = "=(C16<B16)+C16-B16"
$ws.Cells.Item(16,5).Value = 10
$ws.Cells.Item(16,6).Formula = "=(D16*24)*E16"

# New data row 17: 2023-12-27, 15:00 -> 17:00, rate 10
$ws.Cells.Item(17,1).Value = 45287
$ws.Cells.Item(17,2).Value = 0.625
$ws.Cells.Item(17,3).Value = 0.708333333333333
$ws.Cells.Item(17,4).Formula = "=(C17<B17)+C17-B17"
$ws.Cells.Item(17,5).Value = 10
$ws.Cells.Item(17,6).Formula = "=(D17*24)*E17"

# Re-stamp formatting from the template row once more: assigning a Formula to a
# cell can nudge its style (observed on column F), so re-apply the row-15
# formats to guarantee rows 16/17 keep the plain data-row look (not the
# Total-row look) across all six columns.
$ws.Range("A15:F15").Copy()
$ws.Range("A16:F17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fix up the (now relocated) Total row's SUM ranges to cover the two new rows.
$ws.Cells.Item(18,4).Formula = "=SUM(D2:D17)"
$ws.Cells.Item(18,6).Formula = "=SUM(F2:F17)"

# Match the author's final selection (one row below the new Total row).
$ws.Range("F19").Select() | Out-Null

$excel.CutCopyMode = 0
